{"js": "// Update the worksheet date and the twenty-five two-digit multiplication\n// problems to the next day's generated set, per the commit's diff.\nconst replacements = [\n  [\"2024-10-01 Tuesday\", \"2024-10-02 Wednesday\"],\n  [\"68\u00d743=\", \"62\u00d758=\"],\n  [\"76\u00d751=\", \"22\u00d739=\"],\n  [\"57\u00d740=\", \"79\u00d773=\"],\n  [\"86\u00d715=\", \"39\u00d724=\"],\n  [\"22\u00d742=\", \"74\u00d734=\"],\n  [\"38\u00d733=\", \"43\u00d721=\"],\n  [\"77\u00d714=\", \"45\u00d716=\"],\n  [\"30\u00d711=\", \"54\u00d771=\"],\n  [\"55\u00d793=\", \"40\u00d713=\"],\n  [\"74\u00d719=\", \"35\u00d733=\"],\n  [\"25\u00d736=\", \"94\u00d753=\"],\n  [\"35\u00d770=\", \"65\u00d725=\"],\n  [\"23\u00d745=\", \"50\u00d790=\"],\n  [\"33\u00d773=\", \"17\u00d765=\"],\n  [\"61\u00d719=\", \"31\u00d725=\"],\n  [\"53\u00d775=\", \"18\u00d753=\"],\n  [\"26\u00d712=\", \"59\u00d780=\"],\n  [\"82\u00d775=\", \"46\u00d732=\"],\n  [\"98\u00d738=\", \"21\u00d727=\"],\n  [\"32\u00d782=\", \"71\u00d793=\"],\n  [\"88\u00d745=\", \"64\u00d772=\"],\n  [\"93\u00d718=\", \"44\u00d723=\"],\n  [\"74\u00d763=\", \"25\u00d799=\"],\n  [\"89\u00d785=\", \"81\u00d792=\"],\n  [\"51\u00d754=\", \"84\u00d723=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twenty-five two-digit multiplication\n# problems to the next day's generated set, per the commit's diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-01 Tuesday\", \"2024-10-02 Wednesday\"),\n    @(\"68\u00d743=\", \"62\u00d758=\"),\n    @(\"76\u00d751=\", \"22\u00d739=\"),\n    @(\"57\u00d740=\", \"79\u00d773=\"),\n    @(\"86\u00d715=\", \"39\u00d724=\"),\n    @(\"22\u00d742=\", \"74\u00d734=\"),\n    @(\"38\u00d733=\", \"43\u00d721=\"),\n    @(\"77\u00d714=\", \"45\u00d716=\"),\n    @(\"30\u00d711=\", \"54\u00d771=\"),\n    @(\"55\u00d793=\", \"40\u00d713=\"),\n    @(\"74\u00d719=\", \"35\u00d733=\"),\n    @(\"25\u00d736=\", \"94\u00d753=\"),\n    @(\"35\u00d770=\", \"65\u00d725=\"),\n    @(\"23\u00d745=\", \"50\u00d790=\"),\n    @(\"33\u00d773=\", \"17\u00d765=\"),\n    @(\"61\u00d719=\", \"31\u00d725=\"),\n    @(\"53\u00d775=\", \"18\u00d753=\"),\n    @(\"26\u00d712=\", \"59\u00d780=\"),\n    @(\"82\u00d775=\", \"46\u00d732=\"),\n    @(\"98\u00d738=\", \"21\u00d727=\"),\n    @(\"32\u00d782=\", \"71\u00d793=\"),\n    @(\"88\u00d745=\", \"64\u00d772=\"),\n    @(\"93\u00d718=\", \"44\u00d723=\"),\n    @(\"74\u00d763=\", \"25\u00d799=\"),\n    @(\"89\u00d785=\", \"81\u00d792=\"),\n    @(\"51\u00d754=\", \"84\u00d723=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
